$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, pushing existing data down (old rows 2-13 -> 3-14)
$ws.Rows.Item(2).Insert()

# Fill in the new row's data
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Member Engineering"

# Renumber the id column (A) for the rows that got shifted down (now rows 3-14)
for ($i = 3; $i -le 14; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Column B keeps its centered alignment after the insert; make column A match it
# for every data row (id column is now centered just like departmentID)
$ws.Range("A2:A14").HorizontalAlignment = -4108
$ws.Range("B2:B14").HorizontalAlignment = -4108

# Update the selection to match the new active cell
$ws.Range("E13").Select()
